$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (llegadas) - updated statistics
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 98.60000000000002
$ws.Range("D2").Value = 127.0248524160265
$ws.Range("F2").Value = 18.00000000000042
$ws.Range("G2").Value = 57.4999999999978
$ws.Range("H2").Value = 124.9999999999964
$ws.Range("J2").Value = 2.829074033622335
$ws.Range("K2").Value = 10.00076246186234

# Row 3 (servicio) - updated statistics
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 91.63999999999996
$ws.Range("D3").Value = 43.75702425573778
$ws.Range("G3").Value = 83.50000000000026
$ws.Range("H3").Value = 105.9999999999966
$ws.Range("J3").Value = 2.927767384804515
$ws.Range("K3").Value = 12.89274220545711
